# MonthlyLM028.xlsx — 2021/12/09~10 edit by 智偉
# Replace the old "RepayBank / LoanTermMm / LoanTermDd" field block (rows 18-20)
# with the new "MaturityYear / MaturityMonth / MaturityDay" fields, and replace
# the old "BaseRateCode" field (row 26) with a new "ProdNo" field. Both edits
# carry a dated remark in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Row 18: RepayBank/扣款銀行 -> MaturityYear/額度主檔到期日-年 ---
$ws.Range("B18").Value = "MaturityYear"
$ws.Range("C18").Value = "額度主檔到期日-年"
$ws.Range("D18").Value = "DECIMAL"
$ws.Range("E18").Value = 4
$ws.Range("H18").Value = "2021-12-09 智偉修改"

# --- Row 19: LoanTermMm/貸款期間－月 -> MaturityMonth/額度主檔到期日-月 ---
$ws.Range("B19").Value = "MaturityMonth"
$ws.Range("C19").Value = "額度主檔到期日-月"
$ws.Range("D19").Value = "DECIMAL"
$ws.Range("E19").Value = 2
$ws.Range("H19").Value = "2021-12-09 智偉修改"

# --- Row 20: LoanTermDd/貸款期間－日 -> MaturityDay/額度主檔到期日-日 ---
$ws.Range("B20").Value = "MaturityDay"
$ws.Range("C20").Value = "額度主檔到期日-日"
$ws.Range("D20").Value = "DECIMAL"
$ws.Range("E20").Value = 2
$ws.Range("H20").Value = "2021-12-09 智偉修改"

# --- Row 26: BaseRateCode/基本利率代碼 -> ProdNo/商品代碼 ---
$ws.Range("B26").Value = "ProdNo"
$ws.Range("C26").Value = "商品代碼"
$ws.Range("E26").Value = 5
$ws.Range("H26").Value = "2021-12-10 智偉修改"

# --- Column G is no longer a wide free-text remark column; narrow it back down ---
$ws.Columns(7).ColumnWidth = 14.33

# --- Restore the view: scrolled up to row 13, current selection on F27 ---
$ws.Activate()
$ws.Range("F27").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
